$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Comments: bump every "lastCell=J6" -> "lastCell=J4", drop the old
#    per-section comment (A3) and repoint the lignes-formulaire comment (A4)
#    from "section.lignesFormulaire" to "formulaire.lignesFormulaire" (the
#    "section" concept under a formulaire is gone).
# ---------------------------------------------------------------------------
$areaComment = $ws.Range("A1").Comment
$areaComment.Text(("Auteur:" + [char]10 + 'jx:area(lastCell="J4")'))

$formulaireComment = $ws.Range("A2").Comment
$formulaireComment.Text(("Auteur:" + [char]10 + 'jx:each(items="formulaires", var="formulaire", multisheet="sheetNames", lastCell="J4")'))

$sectionComment = $ws.Range("A3").Comment
if ($sectionComment -ne $null) {
    $sectionComment.Delete()
}

$ligneComment = $ws.Range("A4").Comment
$ligneComment.Text(("Auteur:" + [char]10 + 'jx:each(items="formulaire.lignesFormulaire", var="ligne", lastCell="J4")'))

# ---------------------------------------------------------------------------
# 2) Shared strings: the "Section" row is replaced by a "Reponses" banner
#    with Code / Libelle / Reponse column headers.
# ---------------------------------------------------------------------------
$ws.Range("C3:I3").UnMerge()

$ws.Range("A3").Value = "Reponses"
$ws.Range("B3").Value = "Code"
$ws.Range("C3").Value = "Libelle"
$ws.Range("G3").Value = "Reponse"
$ws.Range("D3").Value = $null
$ws.Range("E3").Value = $null
$ws.Range("F3").Value = $null
$ws.Range("H3").Value = $null
$ws.Range("I3").Value = $null

$ws.Range("C3:F3").Merge()
$ws.Range("G3:I3").Merge()

# ---------------------------------------------------------------------------
# 3) Row 1 title banner becomes taller and the label is bold + vertically
#    centred.
# ---------------------------------------------------------------------------
$ws.Rows(1).RowHeight = 39
$ws.Range("A1:J1").Font.Bold = $true
$ws.Range("A1:J1").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 4) "Libelle" / "Numero" field labels on row 2 become bold, and the boxed
#    area spanning rows 2-3 loses the middle divider (row2's bottom border).
# ---------------------------------------------------------------------------
$ws.Range("A2").Font.Bold = $true
$ws.Range("G2").Font.Bold = $true
$ws.Range("B2:F2").Borders(9).LineStyle = 0
$ws.Range("H2:I2").Borders(9).LineStyle = 0

# ---------------------------------------------------------------------------
# 5) New "Reponses" header row (row 3): bold, centred/left, shaded fill and
#    boxed borders, vertically centred - same treatment as the other banner
#    rows.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A3:J3")
$headerRange.Font.Bold = $true
$headerRange.Interior.Pattern = 1
$headerRange.Interior.PatternColorIndex = -4105
$headerRange.Interior.ThemeColor = 5
$headerRange.Interior.TintAndShade = 0
$headerRange.VerticalAlignment = -4108
$ws.Range("A3:B3").HorizontalAlignment = -4131
$ws.Range("C3:I3").HorizontalAlignment = -4131
$ws.Range("A3:J3").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 6) Row 4 (ligne.* formula placeholders) keeps its previous plain styling;
#    only A4 now exists as an (empty) cell carrying the row's base format.
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = $null

# ---------------------------------------------------------------------------
# 7) Sheet cosmetics: column A grows to fit the new "Reponses" banner,
#    the page is set to portrait, and the saved selection moves to H8.
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 10.86
$ws.PageSetup.Orientation = 1
$ws.Range("H8").Select()
